# Atualizado por script em 24-11-2023 20:45
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-MatchRow {
    param(
        [int]$Row,
        [string]$Home,
        [int]$HomeGoals,
        [string]$Away,
        [int]$AwayGoals,
        [double]$HomeOpenOdds,
        [string]$HomeOpenDate,
        [double]$HomeCloseOdds,
        [string]$HomeCloseDate,
        [double]$DrawOpenOdds,
        [string]$DrawOpenDate,
        [double]$DrawCloseOdds,
        [string]$DrawCloseDate,
        [double]$AwayOpenOdds,
        [string]$AwayOpenDate,
        [double]$AwayCloseOdds,
        [string]$AwayCloseDate,
        [string]$Url
    )

    $ws.Cells.Item($Row, 6).Value = $Home
    $ws.Cells.Item($Row, 7).Value = $HomeGoals
    $ws.Cells.Item($Row, 8).Value = $Away
    $ws.Cells.Item($Row, 9).Value = $AwayGoals
    $ws.Cells.Item($Row, 10).Value = $HomeOpenOdds
    $ws.Cells.Item($Row, 11).Value = $HomeOpenDate
    $ws.Cells.Item($Row, 12).Value = $HomeCloseOdds
    $ws.Cells.Item($Row, 13).Value = $HomeCloseDate
    $ws.Cells.Item($Row, 14).Value = $DrawOpenOdds
    $ws.Cells.Item($Row, 15).Value = $DrawOpenDate
    $ws.Cells.Item($Row, 16).Value = $DrawCloseOdds
    $ws.Cells.Item($Row, 17).Value = $DrawCloseDate
    $ws.Cells.Item($Row, 18).Value = $AwayOpenOdds
    $ws.Cells.Item($Row, 19).Value = $AwayOpenDate
    $ws.Cells.Item($Row, 20).Value = $AwayCloseOdds
    $ws.Cells.Item($Row, 21).Value = $AwayCloseDate
    $ws.Cells.Item($Row, 22).Value = $Url
}

# --- Row 42 <-> Row 43 swap (columns F:V only, A:E untouched) ---
Set-MatchRow 42 "Giresunspor" 0 "Manisa FK" 0 `
    3.87 "13/09/2023 22:19" 4.87 "17/09/2023 14:37" `
    3.68 "13/09/2023 22:19" 4.08 "17/09/2023 14:37" `
    1.83 "13/09/2023 22:19" 1.67 "17/09/2023 14:08" `
    "https://www.betexplorer.com/football/turkey/1-lig/giresunspor-manisa-fk/rRSWBTp3/"

Set-MatchRow 43 "Tuzlaspor" 2 "Adanaspor AS" 1 `
    2.73 "10/09/2023 15:12" 2.73 "17/09/2023 14:57" `
    3.27 "10/09/2023 15:12" 3.43 "17/09/2023 14:53" `
    2.51 "10/09/2023 15:12" 2.59 "17/09/2023 14:57" `
    "https://www.betexplorer.com/football/turkey/1-lig/tuzlaspor-adanaspor-as/xUysk6wj/"

# --- Row 50 <-> Row 51 swap (columns F:V only, A:E untouched) ---
Set-MatchRow 50 "Corum" 1 "Sakaryaspor" 0 `
    2.01 "19/09/2023 16:13" 1.97 "23/09/2023 17:59" `
    3.45 "19/09/2023 16:13" 3.4 "23/09/2023 17:59" `
    3.73 "19/09/2023 16:13" 4.08 "23/09/2023 17:59" `
    "https://www.betexplorer.com/football/turkey/1-lig/corum-fk-sakaryaspor/2ov835hq/"

Set-MatchRow 51 "Manisa FK" 0 "Goztepe" 1 `
    2.07 "17/09/2023 18:13" 2.28 "23/09/2023 10:33" `
    3.34 "17/09/2023 18:13" 3.35 "23/09/2023 10:33" `
    3.45 "17/09/2023 18:13" 3.26 "23/09/2023 10:33" `
    "https://www.betexplorer.com/football/turkey/1-lig/manisa-fk-goztepe/j7K2PnF2/"

# --- Row 53 <-> Row 54 swap (columns F:V only, A:E untouched) ---
Set-MatchRow 53 "Sanliurfaspor" 0 "Kocaelispor" 2 `
    2.31 "19/09/2023 16:13" 2.76 "24/09/2023 17:57" `
    3.29 "19/09/2023 16:13" 3.45 "24/09/2023 17:57" `
    3.15 "19/09/2023 16:13" 2.55 "24/09/2023 16:28" `
    "https://www.betexplorer.com/football/turkey/1-lig/sanliurfaspor-kocaelispor/lvbX97FF/"

Set-MatchRow 54 "Adanaspor AS" 1 "Erzurumspor" 0 `
    2.05 "17/09/2023 18:13" 2.64 "24/09/2023 17:57" `
    3.43 "17/09/2023 18:13" 3.28 "24/09/2023 17:59" `
    3.61 "17/09/2023 18:13" 2.78 "24/09/2023 17:57" `
    "https://www.betexplorer.com/football/turkey/1-lig/adanaspor-as-erzurumspor-fk/AmaTAm09/"

# --- Row 96 <-> Row 97 swap (columns F:V only, A:E untouched) ---
Set-MatchRow 96 "Bodrumspor" 2 "Manisa FK" 0 `
    1.95 "29/10/2023 11:42" 1.97 "04/11/2023 20:42" `
    3.43 "29/10/2023 11:42" 3.44 "04/11/2023 20:42" `
    3.98 "29/10/2023 11:42" 4.03 "04/11/2023 20:42" `
    "https://www.betexplorer.com/football/turkey/1-lig/bodrumspor-manisa-fk/WIu9cXQo/"

Set-MatchRow 97 "Bandirmaspor" 2 "Altay" 0 `
    1.53 "29/10/2023 17:13" 1.33 "05/11/2023 11:28" `
    4.25 "29/10/2023 17:13" 5.29 "05/11/2023 11:28" `
    5.78 "29/10/2023 17:13" 9.06 "05/11/2023 11:28" `
    "https://www.betexplorer.com/football/turkey/1-lig/bandirmaspor-altay/4WwHeBeb/"

# --- Row 106 <-> Row 107 swap (columns F:V only, A:E untouched) ---
Set-MatchRow 106 "Boluspor" 2 "Sakaryaspor" 3 `
    2.79 "05/11/2023 14:12" 2.69 "12/11/2023 11:01" `
    3.15 "05/11/2023 14:12" 3.01 "12/11/2023 10:41" `
    2.65 "05/11/2023 14:12" 2.95 "12/11/2023 11:01" `
    "https://www.betexplorer.com/football/turkey/1-lig/boluspor-sakaryaspor/IcM3n7Ip/"

Set-MatchRow 107 "Manisa FK" 0 "Bandirmaspor" 2 `
    2.18 "05/11/2023 11:42" 2.54 "12/11/2023 11:21" `
    3.41 "05/11/2023 11:42" 3.35 "12/11/2023 11:25" `
    3.31 "05/11/2023 11:42" 2.84 "12/11/2023 11:25" `
    "https://www.betexplorer.com/football/turkey/1-lig/manisa-fk-bandirmaspor/Eq666S9G/"

# --- New row 110: copy formatting from row 109, then populate new match data ---
$ws.Range("A109:V109").Copy() | Out-Null
$ws.Range("A110").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Cells.Item(110, 1).Value = 109
$ws.Cells.Item(110, 2).Value = "turkey"
$ws.Cells.Item(110, 3).Value = "1-lig"
$ws.Cells.Item(110, 4).Value = "2023-2024"
$ws.Cells.Item(110, 5).Value = 45254.75

Set-MatchRow 110 "Sakaryaspor" 2 "Erzurumspor" 2 `
    1.71 "16/11/2023 18:43" 1.76 "24/11/2023 17:51" `
    3.72 "16/11/2023 18:43" 3.46 "24/11/2023 17:51" `
    4.91 "16/11/2023 18:43" 5.26 "24/11/2023 17:51" `
    "https://www.betexplorer.com/football/turkey/1-lig/sakaryaspor-erzurumspor-fk/6HNvlKgk/"
